$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Title (row 5) gets the same value as Name (row 4): "OrientationParticuliere"
$ws.Range("B5").Value = "OrientationParticuliere"

# Date (row 8) bumped to the new generation timestamp
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"
